$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 cleanup -----------------------------------------------------
# The "Audited" (column H) cell on row 6 was an empty placeholder cell;
# it is removed entirely (not merely blanked) as part of this edit.
$ws.Range("H6").ClearContents()

# --- Row 7: new "Clinical Labs Form" submission -------------------------
# A brand new row is appended with data for every column except H
# (Audited), which is left unset/blank just like H6 used to be.
$newRow = $ws.Range("A7:I7")

# Force text storage so numeric/date-looking values ("54654",
# "07/20/2020", "1") are kept as literal strings rather than being
# auto-converted to numbers/dates, matching the rest of the sheet.
$newRow.NumberFormat = "@"

$ws.Range("A7").Value = "54654"
$ws.Range("B7").Value = "07/20/2020"
$ws.Range("C7").Value = "1"
$ws.Range("D7").Value = "asdasd"
$ws.Range("E7").Value = "asdasda"
$ws.Range("F7").Value = "asdasd"
$ws.Range("G7").Value = "asdasd"
$ws.Range("I7").Value = "asdasdasd"

# Restore the default style so the new row carries no stray formatting,
# consistent with the rest of the (unstyled) worksheet.
$newRow.Style = "Normal"
